{"js": "// 1. Remove the stray \"_GoBack\" bookmark that currently sits after\n//    \"...rall amount of testing performed\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. \"Code-driven testing\" bullet: \"the results that are returned are\n//    correct\" -> \"the results returned are correct\".\nconst codeDrivenResults = context.document.body.search(\"that are returned\", { matchCase: true });\nawait context.sync();\ncodeDrivenResults.items[0].insertText(\"returned\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. \"The intent of this paper...\" sentence: change \"the GUI testing\n//    framework tools\" to \"multiple GUI testing framework tools\" and\n//    re-insert the \"_GoBack\" bookmark right after the new word \"multiple\".\nconst sentenceResults = context.document.body.search(\n  \"The intent of this paper is to focus on the GUI testing framework tools cu\",\n  { matchCase: true }\n);\nawait context.sync();\n\nconst sentenceRange = sentenceResults.items[0];\nconst words = sentenceRange.split([\" \"], false, false);\nwords.load(\"text\");\nawait context.sync();\n\nconst theWordRange = words.items.find((w) => w.text === \"the \");\n\n// Drop a temporary bookmark right before \"the \" so the preceding text\n// (\"The intent of this paper is to focus on \") stays in its own run once\n// the replacement below happens (otherwise the engine would merge it back\n// together with the replacement text).\ntheWordRange.getRange(\"Start\").insertBookmark(\"TEMP_SPLIT_MARK\");\nawait context.sync();\n\ntheWordRange.insertText(\"multiple \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate \"multiple\" and drop the real bookmark immediately after it\n// (between \"multiple\" and the following \" GUI testing...\" text).\nconst multipleResults = context.document.body.search(\"multiple\", { matchCase: true });\nawait context.sync();\nmultipleResults.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Clean up the temporary helper bookmark.\ncontext.document.deleteBookmark(\"TEMP_SPLIT_MARK\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark that currently sits after\n#    \"...rall amount of testing performed\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. \"Code-driven testing\" bullet: \"the results that are returned are\n#    correct\" -> \"the results returned are correct\".\n$find = $d.Content.Find\n$find.Text = \"that are returned\"\n$find.Replacement.Text = \"returned\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3. \"The intent of this paper...\" sentence: change \"the GUI testing\n#    framework tools\" to \"multiple GUI testing framework tools\" and\n#    re-insert the \"_GoBack\" bookmark right after the new word \"multiple\".\n$full = $d.Content\n$full.Find.Execute(\"the GUI testing framework\") | Out-Null\n$theRange = $d.Range($full.Start, $full.Start + 3)\n\n# Pin both boundaries of the \"the\" range with bookmarks before replacing its\n# text - this keeps the surrounding runs (\"...focus on \" and \" GUI testing\n# framework tools cu\" / \"rrently...\") from being re-merged into a single run\n# by the text assignment below. The right-hand pin is dropped directly under\n# the real bookmark name so it ends up exactly where it belongs.\n$leftPin = $d.Range($theRange.Start, $theRange.Start)\n$d.Bookmarks.Add(\"TEMP_SPLIT_MARK\", $leftPin) | Out-Null\n\n$rightPin = $d.Range($theRange.End, $theRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $rightPin) | Out-Null\n\n$theRange.Text = \"multiple\"\n\nif ($d.Bookmarks.Exists(\"TEMP_SPLIT_MARK\")) {\n    $d.Bookmarks(\"TEMP_SPLIT_MARK\").Delete()\n}\n"}
